$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel's type
# inference get a temporary Text number format so they are stored as
# inline strings (matching the source data), then the format flag is
# cleared back to Normal so no stray style survives on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '56.517.13'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '2.433.74'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  +0.15%  '
Set-TextValue 'D5' '485.27'
$ws.Range('E5').Value = '  +0.26%  '
Set-TextValue 'D6' '151.01'
$ws.Range('E6').Value = '  +6.38%  '
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').Value = '2.436.19'
$ws.Range('E9').Value = '  -2.18%  '
Set-TextValue 'D10' '0.0993'
$ws.Range('E10').Value = '  +1.82%  '
Set-TextValue 'D11' '5.64'
$ws.Range('E11').Value = '  +1.50%  '
Set-TextValue 'D12' '0.332'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').Value = '2.863.51'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '56.677.05'
$ws.Range('E15').Value = '  +1.95%  '
Set-TextValue 'D16' '20.85'
$ws.Range('E16').Value = '  +0.15%  '
Set-TextValue 'D17' '0.0000136'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '2.441.51'
$ws.Range('E18').Value = '  -0.57%  '
Set-TextValue 'D19' '4.55'
$ws.Range('E19').Value = '  +3.48%  '
Set-TextValue 'D20' '322.02'
$ws.Range('E20').Value = '  +1.10%  '
Set-TextValue 'D21' '10.01'
$ws.Range('E21').Value = '  -0.79%  '
Set-TextValue 'D22' '0.997'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  +0.25%  '
Set-TextValue 'D24' '57.72'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D26' '0.403'
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D27' '0.162'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').Value = '2.555.87'
$ws.Range('E28').Value = '  +0.67%  '
Set-TextValue 'D29' '7.41'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('D30').Value = '0.0₃0800'
$ws.Range('E30').Value = '  +2.44%  '
$ws.Range('E31').Value = '  +0.33%  '
Set-TextValue 'D32' '149.80'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('E33').Value = '  +1.79%  '
Set-TextValue 'D34' '18.03'
$ws.Range('E34').Value = '  -1.70%  '
Set-TextValue 'D35' '5.18'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '1.14'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D37' '0.877'
$ws.Range('E37').Value = '  +2.11%  '
Set-TextValue 'D38' '3.70'
$ws.Range('E38').Value = '  +1.18%  '
Set-TextValue 'D39' '1.38'
$ws.Range('E39').Value = '  +4.91%  '
Set-TextValue 'D40' '33.95'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D41' '0.999'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D42' '3.46'
$ws.Range('E42').Value = '  -0.49%  '
Set-TextValue 'D43' '0.0555'
$ws.Range('E43').Value = '  +0.84%  '
Set-TextValue 'D44' '0.603'
$ws.Range('E44').Value = '  -0.73%  '
Set-TextValue 'D45' '0.0954'
$ws.Range('E45').Value = '  +5.90%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D46' '264.53'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D47' '4.80'
$ws.Range('E47').Value = '  +0.50%  '
Set-TextValue 'D48' '10.25'
$ws.Range('E48').Value = '  +0.74%  '
Set-TextValue 'D49' '0.0227'
$ws.Range('E49').Value = '  +0.94%  '
Set-TextValue 'D50' '17.61'
$ws.Range('E50').Value = '  +1.91%  '
Set-TextValue 'D51' '1.72'
$ws.Range('E51').Value = '  +22.24%  '
